$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.090.74'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '1.624.01'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0632'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('E9').Value = '  -1.59%  '
$ws.Range('E10').Value = '  +0.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0849'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').Value = '1.851.44'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.15'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.604.72'
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.88'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.12%  '
$ws.Range('D17').Value = '27.039.95'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '214.30'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.02%  '
$ws.Range('E22').Value = '  -0.95%  '
$ws.Range('E23').Value = '  -6.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('E28').Value = '  -2.73%  '
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.36'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.78%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.746'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +36.32%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '1.356.88'
$ws.Range('E35').Value = '  +3.72%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('E39').Value = '  -1.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.805'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.13'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.29%  '
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('D45').Value = '1.762.25'
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.881'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +32.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.11'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.67%  '
$ws.Range('E48').Value = '  +2.87%  '
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('E50').Value = '  +5.69%  '
$ws.Range('E51').Value = '  +0.41%  '
